$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16 and 17 ("Knärot"/Goodyera repens and "Gultoppig fingersvamp"/
# Ramaria testaceoflava) had their per-observation data swapped: id,
# taxon-sort key, red-list status, taxon id, species/scientific names,
# author, count, unit, easting/northing and start/end time all trade
# places between the two rows, while the columns that already held the
# same value on both rows (validation status, age/stage "L" column
# aside, place name, county, municipality, dates, reporter, ...) are
# left untouched.
#
# Cells are swapped with Cut (rather than read-the-value-then-assign)
# so that cell *type* is preserved exactly as stored (in particular the
# digit-only text in column I, "20"/"25", must stay text and not turn
# into a number), and so no incidental number-format/style gets stamped
# onto the cells.

function Swap-Cell {
  param($addr1, $addr2, $temp)
  $ws.Range($addr1).Cut($ws.Range($temp))
  $ws.Range($addr2).Cut($ws.Range($addr1))
  $ws.Range($temp).Cut($ws.Range($addr2))
}

$tempCell = "ZZ1000"
$cols = @("A","B","D","E","F","G","H","I","J","Q","R","Z","AB")
foreach ($col in $cols) {
  $addr16 = $col + "16"
  $addr17 = $col + "17"
  Swap-Cell $addr16 $addr17 $tempCell
}

# Column L ("Kön") was an empty cell present on row 16 and absent on
# row 17; after the edit it is absent on row 16 and an empty cell on
# row 17 - i.e. it physically relocates from L16 to L17.
$ws.Range("L16").Cut($ws.Range("L17"))
